# Auto stash before merge of "master" and "Robert-Sheehy/master"
#
# Reconstructs the target edit against "Project Documents/Penguin Description.xlsx":
#   - adds a new empty worksheet "Sheet1" at the end of the workbook
#   - activates the first sheet ("Template for Items (2)") and updates its view
#   - widens columns A/B on that sheet
#   - fills in a bunch of new "Can Push / Can Destroy / Can be restrained / ..."
#     rows of the feature table, moving the old "Name of InterFace" row down
#   - applies the built-in "Hyperlink" cell style to the URL cell

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Template for Items (2) -- the first sheet; holds almost all of the edits
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# NOTE: cells are written in the same order the original author entered them
# (not strictly row-major) so that the shared-string table built up the same
# way -- new unique strings land at the same <si> index as the target file.

# -- row 3: role description text is replaced
$ws.Range("B3").Value = "Main Player -Interacts with all Phsical objects in game."

# -- row 6: existing "Eg Turn Left" row gets new copy
$ws.Range("A6").Value = "Can turn 360 degres"
$ws.Range("B6").Value = "Player and camera can look all around the scene"

# -- row 7: brand new row
$ws.Range("A7").Value = "Can move 360 degrees"
$ws.Range("B7").Value = "Player can move 360 degrees through a combination of directional controls"

# -- row 11: existing "Eg Push" row gains full detail
$ws.Range("A11").Value = "Can Push"
$ws.Range("B11").Value = "Player can activate/interact with a certain object(s) to make that object accelerate in a desired direction"
$ws.Range("C11").Value = "Collision"
$ws.Range("E11").Value = "Ice Block"

# -- row 12: brand new row
$ws.Range("A12").Value = "Can Destroy"
$ws.Range("B12").Value = "Player can activate/interact with an object to destroy it"
$ws.Range("C12").Value = "Collision"
$ws.Range("E12").Value = "Pick-up Items, Ice Block"

# -- row 15: existing "ShouldTurnLeft" row gains full detail
$ws.Range("A15").Value = "Can be restrained in movement"
$ws.Range("B15").Value = "Player can be interrupted during movement when certain objects impede the player"
$ws.Range("C15").Value = "Stop Moving"
$ws.Range("E15").Value = "Enemies, World, Ice Block, Rock"

# -- row 16: brand new row (A/B/E filled first, C filled later -- see below)
$ws.Range("A16").Value = "Can be damaged/destroyed"
$ws.Range("B16").Value = "Player can be hurt/killed causing a loss of previously accumulated points or death"
$ws.Range("E16").Value = "Enemies, World."

# -- row 17: brand new "Position / spawn" row (previous row17 content moves to
#    row19 below); A/B filled first, C filled later -- see below
$ws.Range("A17").Value = "Position"
$ws.Range("B17").Value = "Player is spawned on to the map by the world"

# -- back-fill the "Internal Functionality" / "Parameters" column for the two
#    rows just entered
$ws.Range("C16").Value = "(-)H.P. or Death"
$ws.Range("C17").Value = "Poosition"
$ws.Range("E17").Value = "World"

# -- row 19 (row 18 left blank): the old "Name of InterFace" row, now with a
#    different parameters string in B
$ws.Range("A19").Value = "Name of InterFace "
$ws.Range("B19").Value = "Damageable, Killable, Movement, Spawnable, CanPush, CanLook, CanDestroy"

# apply the built-in Hyperlink style to B17 (text only -- no live hyperlink)
$ws.Range("B17").Style = "Hyperlink"

# widen the first two columns
$ws.Columns.Item(1).ColumnWidth = 29.140625
$ws.Columns.Item(2).ColumnWidth = 92.28515625

# make this the active sheet/view, move the selection and reset the zoom
$ws.Activate()
$ws.Range("A20").Select()
$excel.ActiveWindow.Zoom = 100

# ---------------------------------------------------------------------------
# 2. Add the new trailing empty worksheet "Sheet1"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Sheet1"

# re-activate "Template for Items (2)" so it (not the new sheet) is the
# tab shown when the workbook is reopened
$ws.Activate()
